$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 64; this shifts existing rows 64-150 down to 65-151
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new weekly record.
# Static columns mirror the surrounding rows (same market/product).
$ws.Cells.Item(64, 1).Value = 10
$ws.Cells.Item(64, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(64, 3).Value = "La Araucanía"
$ws.Cells.Item(64, 4).Value = 44721
$ws.Cells.Item(64, 5).Value = 9
$ws.Cells.Item(64, 6).Value = 100112012
$ws.Cells.Item(64, 7).Value = "Espinaca"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 30
$ws.Cells.Item(64, 11).Value = 10000
$ws.Cells.Item(64, 12).Value = 10000
$ws.Cells.Item(64, 13).Value = 10000
$ws.Cells.Item(64, 14).Value = "$/docena de atados"
$ws.Cells.Item(64, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(64, 16).Value = 3333
$ws.Cells.Item(64, 17).Value = 3
$ws.Cells.Item(64, 18).Value = "Hortaliza"
